$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new export timestamp
$ws.Name = "IClientBalance-20241211-074902-"

# Update the date column (G) from 45636 (2024-12-10) to 45637 (2024-12-11)
# for all data rows 2 through 274.
for ($r = 2; $r -le 274; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq 45636) {
        $cell.Value = 45637
    }
}
